$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two worker rows that are no longer part of this statement
# (the account statement now only covers "RUBY ROCIO MIRANDA RODRIGUEZ" / period 1810).
# Row 18 first (higher row index) so row 16's index isn't affected by the shift.
$ws.Rows("18").Delete()
$ws.Rows("16").Delete()

# Update the totals now that the data set changed.
$ws.Range("E11").Value = 1067
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
